$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the starting location header (row 2)
$ws.Range("A2").Value = "Starting Location: 700a Lincoln & Barnes"

# New Madigan bike hours: shift schedule, insert 07:00 Lakewood trip at front
$data = New-Object 'object[,]' 196,4
$data[0,0] = "100a"
$data[0,1] = "Madigan Hospital"
$data[0,2] = "to Lakewood"
$data[0,3] = "07:00:00"
$data[1,0] = "100a"
$data[1,1] = "Madigan Hospital"
$data[1,2] = "to Lakewood"
$data[1,3] = "07:00:00"
$data[2,0] = "100a"
$data[2,1] = "Madigan Hospital"
$data[2,2] = "to Lakewood"
$data[2,3] = "07:00:00"
$data[3,0] = "100a"
$data[3,1] = "Madigan Hospital"
$data[3,2] = "to Lakewood"
$data[3,3] = "07:00:00"
$data[4,0] = "801a"
$data[4,1] = "512 Park & Ride"
$data[4,2] = "to Madigan"
$data[4,3] = "07:20:00"
$data[5,0] = "801a"
$data[5,1] = "512 Park & Ride"
$data[5,2] = "to Madigan"
$data[5,3] = "07:20:00"
$data[6,0] = "801a"
$data[6,1] = "512 Park & Ride"
$data[6,2] = "to Madigan"
$data[6,3] = "07:20:00"
$data[7,0] = "801a"
$data[7,1] = "512 Park & Ride"
$data[7,2] = "to Madigan"
$data[7,3] = "07:20:00"
$data[8,0] = "100a"
$data[8,1] = "Madigan Hospital"
$data[8,2] = "to DuPont Transit Center"
$data[8,3] = "07:40:00"
$data[9,0] = "100a"
$data[9,1] = "Madigan Hospital"
$data[9,2] = "to DuPont Transit Center"
$data[9,3] = "07:40:00"
$data[10,0] = "100a"
$data[10,1] = "Madigan Hospital"
$data[10,2] = "to DuPont Transit Center"
$data[10,3] = "07:40:00"
$data[11,0] = "100a"
$data[11,1] = "Madigan Hospital"
$data[11,2] = "to DuPont Transit Center"
$data[11,3] = "07:40:00"
$data[12,0] = "300a"
$data[12,1] = "Post Exchange"
$data[12,2] = "to DuPont Transit Center"
$data[12,3] = "07:50:00"
$data[13,0] = "300a"
$data[13,1] = "Post Exchange"
$data[13,2] = "to DuPont Transit Center"
$data[13,3] = "07:50:00"
$data[14,0] = "300a"
$data[14,1] = "Post Exchange"
$data[14,2] = "to DuPont Transit Center"
$data[14,3] = "07:50:00"
$data[15,0] = "300a"
$data[15,1] = "Post Exchange"
$data[15,2] = "to DuPont Transit Center"
$data[15,3] = "07:50:00"
$data[16,0] = "533a"
$data[16,1] = "DuPont Transit Center"
$data[16,2] = "to Madigan"
$data[16,3] = "08:00:00"
$data[17,0] = "533a"
$data[17,1] = "DuPont Transit Center"
$data[17,2] = "to Madigan"
$data[17,3] = "08:00:00"
$data[18,0] = "533a"
$data[18,1] = "DuPont Transit Center"
$data[18,2] = "to Madigan"
$data[18,3] = "08:00:00"
$data[19,0] = "533a"
$data[19,1] = "DuPont Transit Center"
$data[19,2] = "to Madigan"
$data[19,3] = "08:00:00"
$data[20,0] = "300b"
$data[20,1] = "Post Exchange"
$data[20,2] = "to Madigan"
$data[20,3] = "08:10:00"
$data[21,0] = "300b"
$data[21,1] = "Post Exchange"
$data[21,2] = "to Madigan"
$data[21,3] = "08:10:00"
$data[22,0] = "300b"
$data[22,1] = "Post Exchange"
$data[22,2] = "to Madigan"
$data[22,3] = "08:10:00"
$data[23,0] = "300b"
$data[23,1] = "Post Exchange"
$data[23,2] = "to Madigan"
$data[23,3] = "08:10:00"
$data[24,0] = "100a"
$data[24,1] = "Madigan Hospital"
$data[24,2] = "to McChord Field"
$data[24,3] = "08:20:00"
$data[25,0] = "100a"
$data[25,1] = "Madigan Hospital"
$data[25,2] = "to McChord Field"
$data[25,3] = "08:20:00"
$data[26,0] = "100a"
$data[26,1] = "Madigan Hospital"
$data[26,2] = "to McChord Field"
$data[26,3] = "08:20:00"
$data[27,0] = "100a"
$data[27,1] = "Madigan Hospital"
$data[27,2] = "to McChord Field"
$data[27,3] = "08:20:00"
$data[28,0] = "780a"
$data[28,1] = "Passenger Terminal"
$data[28,2] = "to Madigan"
$data[28,3] = "08:41:20"
$data[29,0] = "780a"
$data[29,1] = "Passenger Terminal"
$data[29,2] = "to Madigan"
$data[29,3] = "08:41:20"
$data[30,0] = "780a"
$data[30,1] = "Passenger Terminal"
$data[30,2] = "to Madigan"
$data[30,3] = "08:41:20"
$data[31,0] = "780a"
$data[31,1] = "Passenger Terminal"
$data[31,2] = "to Madigan"
$data[31,3] = "08:41:20"
$data[32,0] = "100a"
$data[32,1] = "Madigan Hospital"
$data[32,2] = "to DuPont Transit Center"
$data[32,3] = "09:00:00"
$data[33,0] = "100a"
$data[33,1] = "Madigan Hospital"
$data[33,2] = "to DuPont Transit Center"
$data[33,3] = "09:00:00"
$data[34,0] = "100a"
$data[34,1] = "Madigan Hospital"
$data[34,2] = "to DuPont Transit Center"
$data[34,3] = "09:00:00"
$data[35,0] = "100a"
$data[35,1] = "Madigan Hospital"
$data[35,2] = "to DuPont Transit Center"
$data[35,3] = "09:00:00"
$data[36,0] = "300a"
$data[36,1] = "Post Exchange"
$data[36,2] = "to DuPont Transit Center"
$data[36,3] = "09:10:00"
$data[37,0] = "300a"
$data[37,1] = "Post Exchange"
$data[37,2] = "to DuPont Transit Center"
$data[37,3] = "09:10:00"
$data[38,0] = "300a"
$data[38,1] = "Post Exchange"
$data[38,2] = "to DuPont Transit Center"
$data[38,3] = "09:10:00"
$data[39,0] = "300a"
$data[39,1] = "Post Exchange"
$data[39,2] = "to DuPont Transit Center"
$data[39,3] = "09:10:00"
$data[40,0] = "533a"
$data[40,1] = "DuPont Transit Center"
$data[40,2] = "to Madigan"
$data[40,3] = "09:20:00"
$data[41,0] = "533a"
$data[41,1] = "DuPont Transit Center"
$data[41,2] = "to Madigan"
$data[41,3] = "09:20:00"
$data[42,0] = "533a"
$data[42,1] = "DuPont Transit Center"
$data[42,2] = "to Madigan"
$data[42,3] = "09:20:00"
$data[43,0] = "533a"
$data[43,1] = "DuPont Transit Center"
$data[43,2] = "to Madigan"
$data[43,3] = "09:20:00"
$data[44,0] = "300b"
$data[44,1] = "Post Exchange"
$data[44,2] = "to Madigan"
$data[44,3] = "09:30:00"
$data[45,0] = "300b"
$data[45,1] = "Post Exchange"
$data[45,2] = "to Madigan"
$data[45,3] = "09:30:00"
$data[46,0] = "300b"
$data[46,1] = "Post Exchange"
$data[46,2] = "to Madigan"
$data[46,3] = "09:30:00"
$data[47,0] = "300b"
$data[47,1] = "Post Exchange"
$data[47,2] = "to Madigan"
$data[47,3] = "09:30:00"
$data[48,0] = "100a"
$data[48,1] = "Madigan Hospital"
$data[48,2] = "to McChord Field"
$data[48,3] = "09:40:00"
$data[49,0] = "100a"
$data[49,1] = "Madigan Hospital"
$data[49,2] = "to McChord Field"
$data[49,3] = "09:40:00"
$data[50,0] = "100a"
$data[50,1] = "Madigan Hospital"
$data[50,2] = "to McChord Field"
$data[50,3] = "09:40:00"
$data[51,0] = "100a"
$data[51,1] = "Madigan Hospital"
$data[51,2] = "to McChord Field"
$data[51,3] = "09:40:00"
$data[52,0] = "780a"
$data[52,1] = "Passenger Terminal"
$data[52,2] = "to Madigan"
$data[52,3] = "10:01:20"
$data[53,0] = "780a"
$data[53,1] = "Passenger Terminal"
$data[53,2] = "to Madigan"
$data[53,3] = "10:01:20"
$data[54,0] = "780a"
$data[54,1] = "Passenger Terminal"
$data[54,2] = "to Madigan"
$data[54,3] = "10:01:20"
$data[55,0] = "780a"
$data[55,1] = "Passenger Terminal"
$data[55,2] = "to Madigan"
$data[55,3] = "10:01:20"
$data[56,0] = "100a"
$data[56,1] = "Madigan Hospital"
$data[56,2] = "to DuPont Transit Center"
$data[56,3] = "10:20:00"
$data[57,0] = "100a"
$data[57,1] = "Madigan Hospital"
$data[57,2] = "to DuPont Transit Center"
$data[57,3] = "10:20:00"
$data[58,0] = "100a"
$data[58,1] = "Madigan Hospital"
$data[58,2] = "to DuPont Transit Center"
$data[58,3] = "10:20:00"
$data[59,0] = "100a"
$data[59,1] = "Madigan Hospital"
$data[59,2] = "to DuPont Transit Center"
$data[59,3] = "10:20:00"
$data[60,0] = "300a"
$data[60,1] = "Post Exchange"
$data[60,2] = "to DuPont Transit Center"
$data[60,3] = "10:30:00"
$data[61,0] = "300a"
$data[61,1] = "Post Exchange"
$data[61,2] = "to DuPont Transit Center"
$data[61,3] = "10:30:00"
$data[62,0] = "300a"
$data[62,1] = "Post Exchange"
$data[62,2] = "to DuPont Transit Center"
$data[62,3] = "10:30:00"
$data[63,0] = "300a"
$data[63,1] = "Post Exchange"
$data[63,2] = "to DuPont Transit Center"
$data[63,3] = "10:30:00"
$data[64,0] = "533a"
$data[64,1] = "DuPont Transit Center"
$data[64,2] = "to Madigan"
$data[64,3] = "10:40:00"
$data[65,0] = "533a"
$data[65,1] = "DuPont Transit Center"
$data[65,2] = "to Madigan"
$data[65,3] = "10:40:00"
$data[66,0] = "533a"
$data[66,1] = "DuPont Transit Center"
$data[66,2] = "to Madigan"
$data[66,3] = "10:40:00"
$data[67,0] = "533a"
$data[67,1] = "DuPont Transit Center"
$data[67,2] = "to Madigan"
$data[67,3] = "10:40:00"
$data[68,0] = "300b"
$data[68,1] = "Post Exchange"
$data[68,2] = "to Madigan"
$data[68,3] = "10:50:00"
$data[69,0] = "300b"
$data[69,1] = "Post Exchange"
$data[69,2] = "to Madigan"
$data[69,3] = "10:50:00"
$data[70,0] = "300b"
$data[70,1] = "Post Exchange"
$data[70,2] = "to Madigan"
$data[70,3] = "10:50:00"
$data[71,0] = "300b"
$data[71,1] = "Post Exchange"
$data[71,2] = "to Madigan"
$data[71,3] = "10:50:00"
$data[72,0] = "100a"
$data[72,1] = "Madigan Hospital"
$data[72,2] = "to McChord Field"
$data[72,3] = "11:00:00"
$data[73,0] = "100a"
$data[73,1] = "Madigan Hospital"
$data[73,2] = "to McChord Field"
$data[73,3] = "11:00:00"
$data[74,0] = "100a"
$data[74,1] = "Madigan Hospital"
$data[74,2] = "to McChord Field"
$data[74,3] = "11:00:00"
$data[75,0] = "100a"
$data[75,1] = "Madigan Hospital"
$data[75,2] = "to McChord Field"
$data[75,3] = "11:00:00"
$data[76,0] = "780a"
$data[76,1] = "Passenger Terminal"
$data[76,2] = "to Madigan"
$data[76,3] = "11:21:20"
$data[77,0] = "780a"
$data[77,1] = "Passenger Terminal"
$data[77,2] = "to Madigan"
$data[77,3] = "11:21:20"
$data[78,0] = "780a"
$data[78,1] = "Passenger Terminal"
$data[78,2] = "to Madigan"
$data[78,3] = "11:21:20"
$data[79,0] = "780a"
$data[79,1] = "Passenger Terminal"
$data[79,2] = "to Madigan"
$data[79,3] = "11:21:20"
$data[80,0] = "100a"
$data[80,1] = "Madigan Hospital"
$data[80,2] = "to DuPont Transit Center"
$data[80,3] = "11:40:00"
$data[81,0] = "100a"
$data[81,1] = "Madigan Hospital"
$data[81,2] = "to DuPont Transit Center"
$data[81,3] = "11:40:00"
$data[82,0] = "100a"
$data[82,1] = "Madigan Hospital"
$data[82,2] = "to DuPont Transit Center"
$data[82,3] = "11:40:00"
$data[83,0] = "100a"
$data[83,1] = "Madigan Hospital"
$data[83,2] = "to DuPont Transit Center"
$data[83,3] = "11:40:00"
$data[84,0] = "300a"
$data[84,1] = "Post Exchange"
$data[84,2] = "to DuPont Transit Center"
$data[84,3] = "11:50:00"
$data[85,0] = "300a"
$data[85,1] = "Post Exchange"
$data[85,2] = "to DuPont Transit Center"
$data[85,3] = "11:50:00"
$data[86,0] = "300a"
$data[86,1] = "Post Exchange"
$data[86,2] = "to DuPont Transit Center"
$data[86,3] = "11:50:00"
$data[87,0] = "300a"
$data[87,1] = "Post Exchange"
$data[87,2] = "to DuPont Transit Center"
$data[87,3] = "11:50:00"
$data[88,0] = "533a"
$data[88,1] = "DuPont Transit Center"
$data[88,2] = "to Madigan"
$data[88,3] = "12:00:00"
$data[89,0] = "533a"
$data[89,1] = "DuPont Transit Center"
$data[89,2] = "to Madigan"
$data[89,3] = "12:00:00"
$data[90,0] = "533a"
$data[90,1] = "DuPont Transit Center"
$data[90,2] = "to Madigan"
$data[90,3] = "12:00:00"
$data[91,0] = "533a"
$data[91,1] = "DuPont Transit Center"
$data[91,2] = "to Madigan"
$data[91,3] = "12:00:00"
$data[92,0] = "300b"
$data[92,1] = "Post Exchange"
$data[92,2] = "to Madigan"
$data[92,3] = "12:10:00"
$data[93,0] = "300b"
$data[93,1] = "Post Exchange"
$data[93,2] = "to Madigan"
$data[93,3] = "12:10:00"
$data[94,0] = "300b"
$data[94,1] = "Post Exchange"
$data[94,2] = "to Madigan"
$data[94,3] = "12:10:00"
$data[95,0] = "300b"
$data[95,1] = "Post Exchange"
$data[95,2] = "to Madigan"
$data[95,3] = "12:10:00"
$data[96,0] = "100a"
$data[96,1] = "Madigan Hospital"
$data[96,2] = "to McChord Field"
$data[96,3] = "12:20:00"
$data[97,0] = "100a"
$data[97,1] = "Madigan Hospital"
$data[97,2] = "to McChord Field"
$data[97,3] = "12:20:00"
$data[98,0] = "100a"
$data[98,1] = "Madigan Hospital"
$data[98,2] = "to McChord Field"
$data[98,3] = "12:20:00"
$data[99,0] = "100a"
$data[99,1] = "Madigan Hospital"
$data[99,2] = "to McChord Field"
$data[99,3] = "12:20:00"
$data[100,0] = "780a"
$data[100,1] = "Passenger Terminal"
$data[100,2] = "to Madigan"
$data[100,3] = "12:41:20"
$data[101,0] = "780a"
$data[101,1] = "Passenger Terminal"
$data[101,2] = "to Madigan"
$data[101,3] = "12:41:20"
$data[102,0] = "780a"
$data[102,1] = "Passenger Terminal"
$data[102,2] = "to Madigan"
$data[102,3] = "12:41:20"
$data[103,0] = "780a"
$data[103,1] = "Passenger Terminal"
$data[103,2] = "to Madigan"
$data[103,3] = "12:41:20"
$data[104,0] = "100a"
$data[104,1] = "Madigan Hospital"
$data[104,2] = "to DuPont Transit Center"
$data[104,3] = "13:00:00"
$data[105,0] = "100a"
$data[105,1] = "Madigan Hospital"
$data[105,2] = "to DuPont Transit Center"
$data[105,3] = "13:00:00"
$data[106,0] = "100a"
$data[106,1] = "Madigan Hospital"
$data[106,2] = "to DuPont Transit Center"
$data[106,3] = "13:00:00"
$data[107,0] = "100a"
$data[107,1] = "Madigan Hospital"
$data[107,2] = "to DuPont Transit Center"
$data[107,3] = "13:00:00"
$data[108,0] = "300a"
$data[108,1] = "Post Exchange"
$data[108,2] = "to DuPont Transit Center"
$data[108,3] = "13:10:00"
$data[109,0] = "300a"
$data[109,1] = "Post Exchange"
$data[109,2] = "to DuPont Transit Center"
$data[109,3] = "13:10:00"
$data[110,0] = "300a"
$data[110,1] = "Post Exchange"
$data[110,2] = "to DuPont Transit Center"
$data[110,3] = "13:10:00"
$data[111,0] = "300a"
$data[111,1] = "Post Exchange"
$data[111,2] = "to DuPont Transit Center"
$data[111,3] = "13:10:00"
$data[112,0] = "533a"
$data[112,1] = "DuPont Transit Center"
$data[112,2] = "to Madigan"
$data[112,3] = "13:20:00"
$data[113,0] = "533a"
$data[113,1] = "DuPont Transit Center"
$data[113,2] = "to Madigan"
$data[113,3] = "13:20:00"
$data[114,0] = "533a"
$data[114,1] = "DuPont Transit Center"
$data[114,2] = "to Madigan"
$data[114,3] = "13:20:00"
$data[115,0] = "533a"
$data[115,1] = "DuPont Transit Center"
$data[115,2] = "to Madigan"
$data[115,3] = "13:20:00"
$data[116,0] = "300b"
$data[116,1] = "Post Exchange"
$data[116,2] = "to Madigan"
$data[116,3] = "13:30:00"
$data[117,0] = "300b"
$data[117,1] = "Post Exchange"
$data[117,2] = "to Madigan"
$data[117,3] = "13:30:00"
$data[118,0] = "300b"
$data[118,1] = "Post Exchange"
$data[118,2] = "to Madigan"
$data[118,3] = "13:30:00"
$data[119,0] = "300b"
$data[119,1] = "Post Exchange"
$data[119,2] = "to Madigan"
$data[119,3] = "13:30:00"
$data[120,0] = "100a"
$data[120,1] = "Madigan Hospital"
$data[120,2] = "to McChord Field"
$data[120,3] = "13:40:00"
$data[121,0] = "100a"
$data[121,1] = "Madigan Hospital"
$data[121,2] = "to McChord Field"
$data[121,3] = "13:40:00"
$data[122,0] = "100a"
$data[122,1] = "Madigan Hospital"
$data[122,2] = "to McChord Field"
$data[122,3] = "13:40:00"
$data[123,0] = "100a"
$data[123,1] = "Madigan Hospital"
$data[123,2] = "to McChord Field"
$data[123,3] = "13:40:00"
$data[124,0] = "780a"
$data[124,1] = "Passenger Terminal"
$data[124,2] = "to Madigan"
$data[124,3] = "14:01:20"
$data[125,0] = "780a"
$data[125,1] = "Passenger Terminal"
$data[125,2] = "to Madigan"
$data[125,3] = "14:01:20"
$data[126,0] = "780a"
$data[126,1] = "Passenger Terminal"
$data[126,2] = "to Madigan"
$data[126,3] = "14:01:20"
$data[127,0] = "780a"
$data[127,1] = "Passenger Terminal"
$data[127,2] = "to Madigan"
$data[127,3] = "14:01:20"
$data[128,0] = "100a"
$data[128,1] = "Madigan Hospital"
$data[128,2] = "to DuPont Transit Center"
$data[128,3] = "14:20:00"
$data[129,0] = "100a"
$data[129,1] = "Madigan Hospital"
$data[129,2] = "to DuPont Transit Center"
$data[129,3] = "14:20:00"
$data[130,0] = "100a"
$data[130,1] = "Madigan Hospital"
$data[130,2] = "to DuPont Transit Center"
$data[130,3] = "14:20:00"
$data[131,0] = "100a"
$data[131,1] = "Madigan Hospital"
$data[131,2] = "to DuPont Transit Center"
$data[131,3] = "14:20:00"
$data[132,0] = "300a"
$data[132,1] = "Post Exchange"
$data[132,2] = "to DuPont Transit Center"
$data[132,3] = "14:30:00"
$data[133,0] = "300a"
$data[133,1] = "Post Exchange"
$data[133,2] = "to DuPont Transit Center"
$data[133,3] = "14:30:00"
$data[134,0] = "300a"
$data[134,1] = "Post Exchange"
$data[134,2] = "to DuPont Transit Center"
$data[134,3] = "14:30:00"
$data[135,0] = "300a"
$data[135,1] = "Post Exchange"
$data[135,2] = "to DuPont Transit Center"
$data[135,3] = "14:30:00"
$data[136,0] = "533a"
$data[136,1] = "DuPont Transit Center"
$data[136,2] = "to Madigan"
$data[136,3] = "14:40:00"
$data[137,0] = "533a"
$data[137,1] = "DuPont Transit Center"
$data[137,2] = "to Madigan"
$data[137,3] = "14:40:00"
$data[138,0] = "533a"
$data[138,1] = "DuPont Transit Center"
$data[138,2] = "to Madigan"
$data[138,3] = "14:40:00"
$data[139,0] = "533a"
$data[139,1] = "DuPont Transit Center"
$data[139,2] = "to Madigan"
$data[139,3] = "14:40:00"
$data[140,0] = "300b"
$data[140,1] = "Post Exchange"
$data[140,2] = "to Madigan"
$data[140,3] = "14:50:00"
$data[141,0] = "300b"
$data[141,1] = "Post Exchange"
$data[141,2] = "to Madigan"
$data[141,3] = "14:50:00"
$data[142,0] = "300b"
$data[142,1] = "Post Exchange"
$data[142,2] = "to Madigan"
$data[142,3] = "14:50:00"
$data[143,0] = "300b"
$data[143,1] = "Post Exchange"
$data[143,2] = "to Madigan"
$data[143,3] = "14:50:00"
$data[144,0] = "100a"
$data[144,1] = "Madigan Hospital"
$data[144,2] = "to McChord Field"
$data[144,3] = "15:00:00"
$data[145,0] = "100a"
$data[145,1] = "Madigan Hospital"
$data[145,2] = "to McChord Field"
$data[145,3] = "15:00:00"
$data[146,0] = "100a"
$data[146,1] = "Madigan Hospital"
$data[146,2] = "to McChord Field"
$data[146,3] = "15:00:00"
$data[147,0] = "100a"
$data[147,1] = "Madigan Hospital"
$data[147,2] = "to McChord Field"
$data[147,3] = "15:00:00"
$data[148,0] = "780a"
$data[148,1] = "Passenger Terminal"
$data[148,2] = "to Madigan"
$data[148,3] = "15:21:20"
$data[149,0] = "780a"
$data[149,1] = "Passenger Terminal"
$data[149,2] = "to Madigan"
$data[149,3] = "15:21:20"
$data[150,0] = "780a"
$data[150,1] = "Passenger Terminal"
$data[150,2] = "to Madigan"
$data[150,3] = "15:21:20"
$data[151,0] = "780a"
$data[151,1] = "Passenger Terminal"
$data[151,2] = "to Madigan"
$data[151,3] = "15:21:20"
$data[152,0] = "100a"
$data[152,1] = "Madigan Hospital"
$data[152,2] = "to DuPont Transit Center"
$data[152,3] = "15:40:00"
$data[153,0] = "100a"
$data[153,1] = "Madigan Hospital"
$data[153,2] = "to DuPont Transit Center"
$data[153,3] = "15:40:00"
$data[154,0] = "100a"
$data[154,1] = "Madigan Hospital"
$data[154,2] = "to DuPont Transit Center"
$data[154,3] = "15:40:00"
$data[155,0] = "100a"
$data[155,1] = "Madigan Hospital"
$data[155,2] = "to DuPont Transit Center"
$data[155,3] = "15:40:00"
$data[156,0] = "300a"
$data[156,1] = "Post Exchange"
$data[156,2] = "to DuPont Transit Center"
$data[156,3] = "15:50:00"
$data[157,0] = "300a"
$data[157,1] = "Post Exchange"
$data[157,2] = "to DuPont Transit Center"
$data[157,3] = "15:50:00"
$data[158,0] = "300a"
$data[158,1] = "Post Exchange"
$data[158,2] = "to DuPont Transit Center"
$data[158,3] = "15:50:00"
$data[159,0] = "300a"
$data[159,1] = "Post Exchange"
$data[159,2] = "to DuPont Transit Center"
$data[159,3] = "15:50:00"
$data[160,0] = "533a"
$data[160,1] = "DuPont Transit Center"
$data[160,2] = "to Madigan"
$data[160,3] = "16:00:00"
$data[161,0] = "533a"
$data[161,1] = "DuPont Transit Center"
$data[161,2] = "to Madigan"
$data[161,3] = "16:00:00"
$data[162,0] = "533a"
$data[162,1] = "DuPont Transit Center"
$data[162,2] = "to Madigan"
$data[162,3] = "16:00:00"
$data[163,0] = "533a"
$data[163,1] = "DuPont Transit Center"
$data[163,2] = "to Madigan"
$data[163,3] = "16:00:00"
$data[164,0] = "300b"
$data[164,1] = "Post Exchange"
$data[164,2] = "to Madigan"
$data[164,3] = "16:10:00"
$data[165,0] = "300b"
$data[165,1] = "Post Exchange"
$data[165,2] = "to Madigan"
$data[165,3] = "16:10:00"
$data[166,0] = "300b"
$data[166,1] = "Post Exchange"
$data[166,2] = "to Madigan"
$data[166,3] = "16:10:00"
$data[167,0] = "300b"
$data[167,1] = "Post Exchange"
$data[167,2] = "to Madigan"
$data[167,3] = "16:10:00"
$data[168,0] = "100a"
$data[168,1] = "Madigan Hospital"
$data[168,2] = "to Lakewood"
$data[168,3] = "16:20:00"
$data[169,0] = "100a"
$data[169,1] = "Madigan Hospital"
$data[169,2] = "to Lakewood"
$data[169,3] = "16:20:00"
$data[170,0] = "100a"
$data[170,1] = "Madigan Hospital"
$data[170,2] = "to Lakewood"
$data[170,3] = "16:20:00"
$data[171,0] = "100a"
$data[171,1] = "Madigan Hospital"
$data[171,2] = "to Lakewood"
$data[171,3] = "16:20:00"
$data[172,0] = "801a"
$data[172,1] = "512 Park & Ride"
$data[172,2] = "to Madigan"
$data[172,3] = "16:40:00"
$data[173,0] = "801a"
$data[173,1] = "512 Park & Ride"
$data[173,2] = "to Madigan"
$data[173,3] = "16:40:00"
$data[174,0] = "801a"
$data[174,1] = "512 Park & Ride"
$data[174,2] = "to Madigan"
$data[174,3] = "16:40:00"
$data[175,0] = "801a"
$data[175,1] = "512 Park & Ride"
$data[175,2] = "to Madigan"
$data[175,3] = "16:40:00"
$data[176,0] = "100a"
$data[176,1] = "Madigan Hospital"
$data[176,2] = "to DuPont Transit Center"
$data[176,3] = "17:00:00"
$data[177,0] = "100a"
$data[177,1] = "Madigan Hospital"
$data[177,2] = "to DuPont Transit Center"
$data[177,3] = "17:00:00"
$data[178,0] = "100a"
$data[178,1] = "Madigan Hospital"
$data[178,2] = "to DuPont Transit Center"
$data[178,3] = "17:00:00"
$data[179,0] = "100a"
$data[179,1] = "Madigan Hospital"
$data[179,2] = "to DuPont Transit Center"
$data[179,3] = "17:00:00"
$data[180,0] = "300a"
$data[180,1] = "Post Exchange"
$data[180,2] = "to DuPont Transit Center"
$data[180,3] = "17:10:00"
$data[181,0] = "300a"
$data[181,1] = "Post Exchange"
$data[181,2] = "to DuPont Transit Center"
$data[181,3] = "17:10:00"
$data[182,0] = "300a"
$data[182,1] = "Post Exchange"
$data[182,2] = "to DuPont Transit Center"
$data[182,3] = "17:10:00"
$data[183,0] = "300a"
$data[183,1] = "Post Exchange"
$data[183,2] = "to DuPont Transit Center"
$data[183,3] = "17:10:00"
$data[184,0] = "533a"
$data[184,1] = "DuPont Transit Center"
$data[184,2] = "to Madigan"
$data[184,3] = "17:20:00"
$data[185,0] = "533a"
$data[185,1] = "DuPont Transit Center"
$data[185,2] = "to Madigan"
$data[185,3] = "17:20:00"
$data[186,0] = "533a"
$data[186,1] = "DuPont Transit Center"
$data[186,2] = "to Madigan"
$data[186,3] = "17:20:00"
$data[187,0] = "533a"
$data[187,1] = "DuPont Transit Center"
$data[187,2] = "to Madigan"
$data[187,3] = "17:20:00"
$data[188,0] = "300b"
$data[188,1] = "Post Exchange"
$data[188,2] = "to Madigan"
$data[188,3] = "17:30:00"
$data[189,0] = "300b"
$data[189,1] = "Post Exchange"
$data[189,2] = "to Madigan"
$data[189,3] = "17:30:00"
$data[190,0] = "300b"
$data[190,1] = "Post Exchange"
$data[190,2] = "to Madigan"
$data[190,3] = "17:30:00"
$data[191,0] = "300b"
$data[191,1] = "Post Exchange"
$data[191,2] = "to Madigan"
$data[191,3] = "17:30:00"
$data[192,0] = "100a"
$data[192,1] = "Madigan Hospital"
$data[192,2] = "to Lakewood"
$data[192,3] = "17:40:00"
$data[193,0] = "100a"
$data[193,1] = "Madigan Hospital"
$data[193,2] = "to Lakewood"
$data[193,3] = "17:40:00"
$data[194,0] = "100a"
$data[194,1] = "Madigan Hospital"
$data[194,2] = "to Lakewood"
$data[194,3] = "17:40:00"
$data[195,0] = "100a"
$data[195,1] = "Madigan Hospital"
$data[195,2] = "to Lakewood"
$data[195,3] = "17:40:00"

$ws.Range("A4:D199").Value = $data

Write-Host "Applied new Madigan bike hours schedule"
